# LR3/table_1_112.xlsx — reshuffle the summary block at the bottom of the
# sheet: push the four summary rows down by one (leaving row 39 blank),
# append measurement units to their labels, and paint a Times-New-Roman
# "border" block (row 1 and rows 39-47, columns A:K) with the same style
# already used by the data cells (style index 1 == A1's style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture old summary formulas/values before we clear them --------
$sumItogoFormula  = $ws.Range("C43").Formula
$avgAreaFormula   = $ws.Range("C44").Formula
$maxDelayFormula  = $ws.Range("C45").Formula
$maxPaymentFormula = $ws.Range("C46").Formula

# --- 2. Clear the old (pre-move) summary rows 43:46 ----------------------
$ws.Range("A43:K46").ClearContents()

# --- 3. Re-create the four summary rows one row higher up (40:43) with
#        updated (unit-suffixed) labels -----------------------------------
$ws.Range("B40").Value = "общая сумма графы ""Итого"", руб."
$ws.Range("C40").Formula = $sumItogoFormula

$ws.Range("B41").Value = "средняя площадь, кв.м."
$ws.Range("C41").Formula = $avgAreaFormula

$ws.Range("B42").Value = "максимальный срок просрочки, дней"
$ws.Range("C42").Formula = $maxDelayFormula

$ws.Range("B43").Value = "максимальная сумма к оплате, руб."
$ws.Range("C43").Formula = $maxPaymentFormula

# --- 4. Paint the Times-New-Roman style (same as A1 / the data block)
#        across row 1 and rows 39-47, columns A:K --------------------------
$ws.Range("A1").Copy()
$ws.Range("B1:K1").PasteSpecial(-4122)
$ws.Range("A39:K43").PasteSpecial(-4122)
$ws.Range("A44:A46").PasteSpecial(-4122)
$ws.Range("D44:K46").PasteSpecial(-4122)
$ws.Range("A47:K47").PasteSpecial(-4122)

# give the newly styled rows the same (bigger-font) row height as the rest
# of the sheet
$ws.Range("A39:K47").RowHeight = 15.75

# --- 5. Update the sheet view (scrolled down to the summary block, with
#        D42 selected) -----------------------------------------------------
$ws.Range("D42").Select()
$excel.ActiveWindow.ScrollRow = 35
